# Fruta / hortaliza, semanal
#
# The weekly refresh rotates the weekly price-report rows for this
# market/product subset: the data that used to live in rows 2, 7, 3, 8, 6, 9
# shifts forward by one slot in that cyclic order (row 2 -> row 7,
# row 7 -> row 3, row 3 -> row 8, row 8 -> row 6, row 6 -> row 9,
# row 9 -> row 2). Columns A,B,C,E,F,G,H,I,J,K,R are identical for every
# row in this block, so only the varying columns (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Precio $/Kg,
# Kg/unidad) actually need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($ws, $row) {
    [ordered]@{
        D = $ws.Cells.Item($row, 4).Value2
        L = $ws.Cells.Item($row, 12).Value2
        M = $ws.Cells.Item($row, 13).Value2
        N = $ws.Cells.Item($row, 14).Value2
        O = $ws.Cells.Item($row, 15).Value2
        P = $ws.Cells.Item($row, 16).Value2
        Q = $ws.Cells.Item($row, 17).Value2
        S = $ws.Cells.Item($row, 19).Value2
        T = $ws.Cells.Item($row, 20).Value2
    }
}

function Set-RowData($ws, $row, $data) {
    $ws.Cells.Item($row, 4).Value = $data.D
    $ws.Cells.Item($row, 12).Value = $data.L
    $ws.Cells.Item($row, 13).Value = $data.M
    $ws.Cells.Item($row, 14).Value = $data.N
    $ws.Cells.Item($row, 15).Value = $data.O
    $ws.Cells.Item($row, 16).Value = $data.P
    $ws.Cells.Item($row, 17).Value = $data.Q
    $ws.Cells.Item($row, 19).Value = $data.S
    $ws.Cells.Item($row, 20).Value = $data.T
}

# Snapshot the "before" state of every row in the rotation first, since
# several of the destinations also act as sources.
$row2 = Get-RowData $ws 2
$row3 = Get-RowData $ws 3
$row6 = Get-RowData $ws 6
$row7 = Get-RowData $ws 7
$row8 = Get-RowData $ws 8
$row9 = Get-RowData $ws 9

# Apply the cyclic rotation: row N receives the data that used to sit in
# its predecessor in the cycle (2 <- 9 <- 6 <- 8 <- 3 <- 7 <- 2).
Set-RowData $ws 2 $row9
Set-RowData $ws 3 $row7
Set-RowData $ws 6 $row8
Set-RowData $ws 7 $row2
Set-RowData $ws 8 $row3
Set-RowData $ws 9 $row6
